$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift all timestamps in column A (rows 2-97) forward by 18 days.
for ($r = 2; $r -le 97; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $cell.Value = $cell.Value2 + 18
}

# Update the production values (column B) for the new day's solar ramp.
$newValues = @{
    23 = 4
    24 = 23
    25 = 49
    26 = 107
    27 = 172
    28 = 249
    29 = 348
    30 = 513
    31 = 619
    32 = 717
    33 = 814
    34 = 984
    35 = 1095
    36 = 1157
    37 = 1169
    38 = 1294
    39 = 1411
}

foreach ($r in $newValues.Keys) {
    $ws.Cells.Item($r, 2).Value = $newValues[$r]
}
